$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (legacy password hash D382); unprotect so the
# cells below (which are locked, the default) can be edited, then restore
# protection afterwards.
$ws.Unprotect()

# Bump the "as of" date in the confidential disclaimer footer (A10):
# 2021-04-30 -> 2021-05-03
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-03 for illustrative purposes only and are subject to change."
# Setting multi-line text auto-expands the row height; put it back to the
# sheet's default (the row carried no explicit height before the edit).
$ws.Rows(10).AutoFit()

# Refresh the Weight (D) and Percent Change (E) figures for each fund row.
$ws.Range("D2").Value = 0.2521386453723293
$ws.Range("E2").Value = -0.007490636704119646

$ws.Range("D3").Value = 0.490269459288686
$ws.Range("E3").Value = 0.009641135511515797

$ws.Range("D4").Value = 0.1000774751929986
$ws.Range("E4").Value = -0.008542141230068356

$ws.Range("D5").Value = 0.1003890078208202
$ws.Range("E5").Value = 0.007717750826901959

$ws.Range("D6").Value = 0.05712541232516588
$ws.Range("E6").Value = 0.005155794664873481

$ws.Range("D7").Value = 0.9999999999999999
$ws.Range("E7").Value = 0.003052503619775271

# Restore sheet protection.
$ws.Protect()
